$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -20.517
$ws.Range("C4").Value = -12.88

$ws.Range("A7").Value = -19.796

$ws.Range("D10").Value = -7.675

$ws.Range("C12").Value = -11.138

$ws.Range("D13").Value = -8.061

$ws.Range("A16").Value = -22.059

$ws.Range("C18").Value = -12.612

$ws.Range("C19").Value = -11.77

$ws.Range("C20").Value = -12.067

$ws.Range("A29").Value = -21.32399999999999

$ws.Range("D30").Value = -7.139

$ws.Range("C31").Value = -13.298

$ws.Range("A32").Value = -21.781

$ws.Range("A40").Value = -19.985
$ws.Range("C40").Value = -12.038
$ws.Range("D40").Value = -7.996

$ws.Range("C42").Value = -12.355

$ws.Range("D44").Value = -7.502

$ws.Range("C47").Value = -11.844

$ws.Range("C48").Value = -11.97

$ws.Range("A52").Value = -21.918

$ws.Range("A57").Value = -22.241

$ws.Range("C63").Value = -11.207

$ws.Range("C64").Value = -10.743

$ws.Range("A66").Value = -21.652

$ws.Range("C76").Value = -12.969

$ws.Range("C81").Value = -13.1

$ws.Range("C89").Value = -12.937
$ws.Range("D89").Value = -8.122999999999999

$ws.Range("D91").Value = -7.053

$ws.Range("C94").Value = -11.375

$ws.Range("A100").Value = -22.26

$wb.Save()
